# Add two new columns, I ("I0") and J ("IF"), to the data sheet, matching
# the style of the existing header row and filling in values for every
# data row (2-28).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header style used by the rest of row 1 (bold, centered,
# bordered) by copying the formatting from the existing "IP" header cell.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Data rows (2-28): new I0/IF values ---
$data = @{
    2  = @(1, 4)
    3  = @(1, 5)
    4  = @(1, 5)
    5  = @(1, 5)
    6  = @(2, 5)
    7  = @(9, 9)
    8  = @(6, 8)
    9  = @(7, 7)
    10 = @(7, 7)
    11 = @(6, 7)
    12 = @(7, 8)
    13 = @(7, 8)
    14 = @(7, 9)
    15 = @(7, 8)
    16 = @(4, 5)
    17 = @(5, 7)
    18 = @(6, 7)
    19 = @(8, 8)
    20 = @(5, 5)
    21 = @(7, 7)
    22 = @(4, 6)
    23 = @(7, 7)
    24 = @(7, 7)
    25 = @(1, 3)
    26 = @(8, 9)
    27 = @(6, 7)
    28 = @(8, 8)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}

Write-Output "I0/IF columns added"
